# Edit script: rename sheet, rename header, update/add expense rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab name) from "Expense" to "expense".
$ws.Name = "expense"

# 2. Rename header A1 from "Source" to "category".
$ws.Range("A1").Value = "category"

# 3. Update the existing data row (row 2): Food/250 -> Fun/10000, new date.
$ws.Range("A2").Value = "Fun"
$ws.Range("B2").Value = 10000
$ws.Range("C2").Value = 45801.22928240741

# 4. Add two new data rows (3 and 4), copying the date format from C2 so the
#    numeric date values keep the same display style (m/d/yyyy).
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3:C4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A3").Value = "Books"
$ws.Range("B3").Value = 500
$ws.Range("C3").Value = 45778.22928240741

$ws.Range("A4").Value = "food"
$ws.Range("B4").Value = 300
$ws.Range("C4").Value = 45717.22928240741
